$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep numeric-looking "Price" strings as text (matches the source inlineStr typing)
# instead of letting Excel auto-convert them to numbers.
$ws.Range('D2:D40').NumberFormat = '@'
$ws.Range('D42:D51').NumberFormat = '@'

$ws.Range('D2').Value = '28.260.77'
$ws.Range('E2').Value = '  +3.02%  '
$ws.Range('D3').Value = '1.834.16'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('D4').Value = '0.9975'
$ws.Range('E4').Value = '  -0.50%  '
$ws.Range('D5').Value = '340.21'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').Value = '0.9946'
$ws.Range('D7').Value = '0.3947'
$ws.Range('E7').Value = '  +3.62%  '
$ws.Range('D8').Value = '0.3509'
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('D9').Value = '48.30'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').Value = '1.207'
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').Value = '0.07633'
$ws.Range('E11').Value = '  +1.36%  '
$ws.Range('D12').Value = '0.9960'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = '22.35'
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('D14').Value = '6.584'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').Value = '1.825.41'
$ws.Range('D16').Value = '7.259'
$ws.Range('E16').Value = '  +2.59%  '
$ws.Range('D17').Value = '0.00001114'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '0.06710'
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').Value = '85.93'
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').Value = '0.9963'
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').Value = '18.05'
$ws.Range('E21').Value = '  +3.76%  '
$ws.Range('D22').Value = '6.629'
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').Value = '28.224.44'
$ws.Range('E23').Value = '  +2.99%  '
$ws.Range('D24').Value = '12.85'
$ws.Range('E24').Value = '  +2.23%  '
$ws.Range('D25').Value = '2.400'
$ws.Range('E25').Value = '  -1.61%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '2.616'
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').Value = '1.526'
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('D28').Value = '21.66'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('D29').Value = '155.48'
$ws.Range('E29').Value = '  +1.63%  '
$ws.Range('D30').Value = '2.031.68'
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('D31').Value = '136.62'
$ws.Range('E31').Value = '  +1.60%  '
$ws.Range('D32').Value = '6.302'
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('D33').Value = '4.037'
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('D34').Value = '0.08884'
$ws.Range('E34').Value = '  +1.92%  '
$ws.Range('D35').Value = '13.40'
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('D36').Value = '5.602'
$ws.Range('E36').Value = '  +2.17%  '
$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').Value = '0.7033'
$ws.Range('E37').Value = '  +1.67%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.02455'
$ws.Range('E38').Value = '  +4.84%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.06609'
$ws.Range('E39').Value = '  +3.36%  '
$ws.Range('D40').Value = '1.613'
$ws.Range('E40').Value = '  -4.83%  '
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('D42').Value = '1.274'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').Value = '8.657'
$ws.Range('E43').Value = '  -3.44%  '
$ws.Range('D44').Value = '14.66'
$ws.Range('E44').Value = '  +1.69%  '
$ws.Range('D45').Value = '0.6560'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D46').Value = '3.894'
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('D47').Value = '2.187'
$ws.Range('E47').Value = '  +2.38%  '
$ws.Range('D48').Value = '132.82'
$ws.Range('E48').Value = '  +1.97%  '
$ws.Range('D49').Value = '0.07247'
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('D50').Value = '81.00'
$ws.Range('E50').Value = '  +1.73%  '
$ws.Range('D51').Value = '1.169'
$ws.Range('E51').Value = '  +3.90%  '
